# Update automatico via Actualizar 02-05-2021 19-27-09
#
# The sheet holds a repeating 14-row block (one row per monitored
# service). This run:
#   1. refreshes the timestamp (col D) of the most recent block
#      (rows 744-757) to the new check time, and
#   2. appends a brand-new block (rows 758-771) for the next check,
#      cloned from the previous block (names/URLs/"Disponible" are
#      identical - only the timestamp changes) including the B-column
#      hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstSrcRow  = 744
$lastSrcRow   = 757
$blockSize    = $lastSrcRow - $firstSrcRow + 1   # 14
$firstNewRow  = 758

$refreshedStamp = 44232.78916576389
$newStamp       = 44232.81028977159

# Target URL (and optional #fragment "location") for each of the 14
# cyclical rows, in row order, matching the existing hyperlink targets.
$targets = @(
  @{ Url = "https://www.dataintelligence-group.com/";                         Loc = $null },
  @{ Url = "https://serviciodashboard.azurewebsites.net/";                    Loc = $null },
  @{ Url = "https://powerbi.microsoft.com/es-es/";                            Loc = $null },
  @{ Url = "https://www.dropbox.com/";                                        Loc = $null },
  @{ Url = "https://dataintelligence.store/";                                 Loc = $null },
  @{ Url = "https://app-data-i.users.earthengine.app/";                       Loc = $null },
  @{ Url = "https://odooutil.azurewebsites.net/";                             Loc = $null },
  @{ Url = "https://filtradordashboard.azurewebsites.net/";                   Loc = $null },
  @{ Url = "https://ide.dataintelligence-group.com/mapstore/";                Loc = "/" },
  @{ Url = "https://ide.dataintelligence-group.com/geoserver/web/?0";         Loc = $null },
  @{ Url = "https://ide.dataintelligence-group.com/";                         Loc = $null },
  @{ Url = "https://rpubs.com/dataintelligence/";                             Loc = $null },
  @{ Url = "https://github.com/Sud-Austral/";                                 Loc = $null },
  @{ Url = "https://ezexporter.highviewapps.com/exports/export-profile/";     Loc = $null }
)

# 1) Bump the timestamp on the existing last block (744-757).
for ($r = $firstSrcRow; $r -le $lastSrcRow; $r++) {
  $ws.Cells.Item($r, 4).Value = $refreshedStamp
}

# 2) Clone the block into the new rows (758-771): copy the whole row
#    (values + styles) so names/URLs/"Disponible" text and formatting
#    (incl. the Hyperlink style on column B, date format on column D)
#    come along for free, then set the new timestamp.
for ($i = 0; $i -lt $blockSize; $i++) {
  $srcRow = $firstSrcRow + $i
  $dstRow = $firstNewRow + $i
  $ws.Range("A" + $srcRow + ":D" + $srcRow).Copy($ws.Range("A" + $dstRow + ":D" + $dstRow))
  $ws.Cells.Item($dstRow, 4).Value = $newStamp
}

# 3) Re-create the hyperlinks on the new column-B cells.
for ($i = 0; $i -lt $blockSize; $i++) {
  $dstRow = $firstNewRow + $i
  $t = $targets[$i]
  $cell = $ws.Cells.Item($dstRow, 2)
  if ($t.Loc) {
    $ws.Hyperlinks.Add($cell, $t.Url, $t.Loc)
  } else {
    $ws.Hyperlinks.Add($cell, $t.Url)
  }
  # Hyperlinks.Add re-stamps the cell's style; put the original
  # "Hyperlink" cell style (shared with the rest of column B) back.
  $ws.Range("B" + $dstRow).Style = "Hyperlink"
}
